$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename existing species rows to match the author's updated ideas.
$ws.Range("A23").Value = "Toxic"
$ws.Range("A28").Value = "Acidic"

# Append a brand-new "Fungal" ants row (row 29) with its stats.
$ws.Range("A29").Value = "Fungal"
$ws.Range("B29").Value = $false
$ws.Range("C29").Value = 2
$ws.Range("D29").Value = 15
$ws.Range("E29").Value = "Swampland, Jungle"
$ws.Range("F29").Value = $false
$ws.Range("G29").Value = $true
$ws.Range("H29").Value = "Cultivator + Jungle"
$ws.Range("I29").Value = "Grows Fungi"

# Reflect the author's final cursor position in the saved view state.
[void]$ws.Range("I23").Select()
